$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct/"mis-spell" the keyword step text in column C (rows 2-9) in place.
$ws.Range("C2").Value = 'Given  While 1 filling the form, navigation for "Testzen Labs Form" to proceed with registration.'
$ws.Range("C3").Value = 'When you correctly entered   the "First Name" before moving to the next field.'
$ws.Range("C4").Value = 'And You should carefully  entterred the "Last Name" so that it matches your official documents.'
$ws.Range("C5").Value = 'And Before proceeding further, make sure to entterring the "Phone Number" to receive OTP verification.'
$ws.Range("C6").Value = 'Then In the form, selection "Country" from the dropdown list to specify your nationality.'
$ws.Range("C7").Value = 'And To complete your application, kindly upload "Resume" in the specified format.'
$ws.Range("C8").Value = 'Then checked the "Male" option if applicable.'
$ws.Range("C9").Value = 'And generation a random number for the pin code before submission.'

# Move the active selection from C2 to C4.
$ws.Range("C4").Select()
